# The upstream change (commit "Moving from 2.0.2 to 2.0.3") is a pure
# OOXML re-serialization: the document was re-saved by a newer docx4j
# version whose JAXB marshaller emits element attributes (and the root
# namespace declarations) in alphabetical order instead of declaration
# order. No paragraph text, run formatting, styles, headers, footers,
# footnotes, section properties or any other document content actually
# changed between the two revisions -- every hunk in the diff is a
# like-for-like attribute/namespace reordering (e.g.
# `w:headerReference w:type="even" r:id="rId6"` -> the same element with
# `r:id="rId6" w:type="even"`, `w:pgSz w:w=".." w:h=".."` -> `w:h=".." w:w=".."`,
# etc.), so there is nothing in the Word object model (paragraphs, runs,
# sections, styles, headers/footers...) for a COM script to change.
#
# Word's COM automation surface edits document *content*; it has no way
# to dictate the raw XML attribute-serialization order that the
# underlying package writer chooses, so the faithful reproduction of
# this particular diff is a script that touches nothing and leaves the
# document exactly as authored.

$d = $word.ActiveDocument
